$d = $word.ActiveDocument

function Replace-Text($old, $new) {
    $d.Content.Find.Execute($old, $true, $false, $false, $false, $false, `
                             $true, 1, $false, $new, 2)
}

Replace-Text "2025-09-28 Sunday" "2025-09-29 Monday"

Replace-Text "71×29=" "77×53="
Replace-Text "59×24=" "41×27="
Replace-Text "33×92=" "57×20="
Replace-Text "41×21=" "74×52="
Replace-Text "52×59=" "93×17="
Replace-Text "41×16=" "63×73="
Replace-Text "91×18=" "52×51="
Replace-Text "61×72=" "68×71="
Replace-Text "50×65=" "18×85="
Replace-Text "96×33=" "45×85="
Replace-Text "55×25=" "23×69="
Replace-Text "93×36=" "78×72="
Replace-Text "54×66=" "30×22="
Replace-Text "98×11=" "55×29="
Replace-Text "18×52=" "22×70="
Replace-Text "77×96=" "60×80="
Replace-Text "36×63=" "20×99="
Replace-Text "26×74=" "36×72="
Replace-Text "17×93=" "86×50="
Replace-Text "19×72=" "90×83="
Replace-Text "30×24=" "69×89="
Replace-Text "97×83=" "37×72="
Replace-Text "40×15=" "71×80="
Replace-Text "41×22=" "94×78="
Replace-Text "33×68=" "16×18="
